# Batch processing update for study 1: add a "no_cutoff" column (F) that
# records a "none" value for every participant/trial row, reflecting the
# trials that were processed without a cutoff frequency.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F, matching the header style already used in row 1.
$ws.Range("F1").Value = "no_cutoff"

# Data rows 2 through 126 all get the literal value "none" in column F.
$ws.Range("F2:F126").Value = "none"

# The author's selection moved from D7 to H7 after adding the new column.
$ws.Range("H7").Select()
